$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

function Set-TextValue($cellAddr, $val) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "26.946.75"
$ws.Range("E2").Value = "  +1.09%  "

$ws.Range("D3").Value = "1.820.58"
$ws.Range("E3").Value = "  +1.90%  "

$ws.Range("E4").Value = "  -0.23%  "

Set-TextValue "D5" "309.90"
$ws.Range("E5").Value = "  +0.57%  "

Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  -0.25%  "

Set-TextValue "D7" "0.4682"
$ws.Range("E7").Value = "  +3.43%  "

Set-TextValue "D8" "0.3697"
$ws.Range("E8").Value = "  +0.36%  "

Set-TextValue "D9" "0.07381"
$ws.Range("E9").Value = "  +1.74%  "

Set-TextValue "D10" "0.8724"
$ws.Range("E10").Value = "  +2.49%  "

$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("D12").Value = "1.844.07"
$ws.Range("E12").Value = "  +3.13%  "

$ws.Range("E13").Value = "  +1.66%  "

Set-TextValue "D14" "92.74"
$ws.Range("E14").Value = "  +2.76%  "

Set-TextValue "D15" "0.07079"
$ws.Range("E15").Value = "  +0.64%  "

Set-TextValue "D16" "6.518"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("E17").Value = "  -0.21%  "

Set-TextValue "D18" "0.000008727"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("E19").Value = "  -0.35%  "

Set-TextValue "D20" "14.78"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("D21").Value = "26.969.48"
$ws.Range("E21").Value = "  +1.15%  "

Set-TextValue "D22" "5.336"
$ws.Range("E22").Value = "  +1.77%  "

Set-TextValue "D23" "10.58"
$ws.Range("E23").Value = "  -0.69%  "

$ws.Range("D24").Value = "2.044.09"
$ws.Range("E24").Value = "  +1.32%  "

Set-TextValue "D25" "1.905"
$ws.Range("E25").Value = "  -0.03%  "

Set-TextValue "D26" "151.70"
$ws.Range("E26").Value = "  +0.91%  "

Set-TextValue "D27" "2.199"
$ws.Range("E27").Value = "  +1.60%  "

Set-TextValue "D28" "18.46"
$ws.Range("E28").Value = "  +1.69%  "

Set-TextValue "D29" "5.333"
$ws.Range("E29").Value = "  +3.04%  "

Set-TextValue "D30" "116.01"
$ws.Range("E30").Value = "  +1.77%  "

Set-TextValue "D31" "0.08939"
$ws.Range("E31").Value = "  +1.33%  "

Set-TextValue "D32" "0.7714"
$ws.Range("E32").Value = "  +2.17%  "

Set-TextValue "D33" "1.168"
$ws.Range("E33").Value = "  +1.60%  "

Set-TextValue "D34" "4.509"
$ws.Range("E34").Value = "  +1.73%  "

$ws.Range("E35").Value = "  +1.14%  "

Set-TextValue "D36" "1.001"
$ws.Range("E36").Value = "  -0.21%  "

Set-TextValue "D37" "1.087"
$ws.Range("E37").Value = "  -1.92%  "

$ws.Range("E38").Value = "  +1.57%  "

$ws.Range("E39").Value = "  +2.13%  "

Set-TextValue "D40" "7.337"
$ws.Range("E40").Value = "  +3.17%  "

Set-TextValue "D41" "2.950"
$ws.Range("E41").Value = "  +3.27%  "

Set-TextValue "D42" "0.5363"
$ws.Range("E42").Value = "  +3.43%  "

Set-TextValue "D43" "2.373"
$ws.Range("E43").Value = "  +2.64%  "

Set-TextValue "D44" "0.1673"
$ws.Range("E44").Value = "  +1.81%  "

Set-TextValue "D45" "8.468"
$ws.Range("E45").Value = "  +0.39%  "

Set-TextValue "D46" "0.4976"
$ws.Range("E46").Value = "  +0.97%  "

Set-TextValue "D47" "10.50"
$ws.Range("E47").Value = "  +3.20%  "

Set-TextValue "D48" "1.675"
$ws.Range("E48").Value = "  +2.30%  "

$ws.Range("E49").Value = "  -0.30%  "

Set-TextValue "D50" "103.37"
$ws.Range("E50").Value = "  +0.20%  "

Set-TextValue "D51" "0.06295"
$ws.Range("E51").Value = "  +0.38%  "

$scratch.Clear()